$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new expense row (row 4) below the existing data:
# Category=Clothing, Name=wads, Date=2023-03-20, Price=453.0, Account=Checkings
$ws.Range("A4").Value = "Clothing"
$ws.Range("B4").Value = "wads"
# Date/Price look numeric - prefix with a quote so Excel stores them as
# literal text (matching the rest of the sheet) instead of converting
# them to a date serial / number.
$ws.Range("C4").Value = "'2023-03-20"
$ws.Range("D4").Value = "'453.0"
$ws.Range("E4").Value = "Checkings"
